$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18, pushing the existing rows 18-20 down to 19-21
$ws.Rows("18:18").Insert()

# Populate the newly inserted row 18 with the new weekly price record
$ws.Range("A18").Value2 = 1
$ws.Range("B18").Value2 = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C18").Value2 = 'Arica y Parinacota'
$ws.Range("D18").Value2 = 44474
$ws.Range("E18").Value2 = 15
$ws.Range("F18").Value2 = 100112052
$ws.Range("G18").Value2 = 'Albahaca'
$ws.Range("H18").Value2 = 'Sin especificar'
$ws.Range("I18").Value2 = 'Primera'
$ws.Range("J18").Value2 = 250
$ws.Range("K18").Value2 = 2000
$ws.Range("L18").Value2 = 2500
$ws.Range("M18").Value2 = 2250
$ws.Range("N18").Value2 = '$/paquete'
$ws.Range("O18").Value2 = 'Región de Arica y Parinacota'
$ws.Range("P18").Value2 = 2250
$ws.Range("Q18").Value2 = 1
$ws.Range("R18").Value2 = 'Hortaliza'
